$d = $word.ActiveDocument

# Constant used for red font color (wdColorRed = 255 -> 0x0000FF, i.e. RGB(255,0,0))
$RED = 255

# ---------------------------------------------------------------------------
# Edit 1: "Conducting poison test every time " -> "Conducting a poison test every time "
#         ("a" inserted in red)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Conducting ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Collapse(0)  # wdCollapseEnd -> right after "Conducting "
    $insertStart1 = $r1.Start
    $r1.InsertAfter("a ")
    $aRange = $d.Range($insertStart1, $insertStart1 + 1)
    $aRange.Font.Color = $RED
}

# ---------------------------------------------------------------------------
# Edit 2: "looking its color, shape, habitat etc.," -> "looking at its color, shape, habitat etc.,"
#         ("at" inserted in red)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("looking", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)  # right after "looking"
    $insertStart2 = $r2.Start
    $r2.InsertAfter(" at")
    $atRange = $d.Range($insertStart2 + 1, $insertStart2 + 3)
    $atRange.Font.Color = $RED
}

# ---------------------------------------------------------------------------
# Edit 3: "is really required" -> "is really essential" ("required" replaced
#         by "essential" in red); the _GoBack bookmark keeps marking the
#         position right before the replaced word instead of right after it.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("required", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $reqStart = $r3.Start
    $r3.Delete()

    $bm = $d.Bookmarks.Item("_GoBack")
    $bmStart = $bm.Start
    $bm.Delete()

    $insertPoint = $d.Range($bmStart, $bmStart)
    $insertPoint.InsertAfter("essential")
    $essentialRange = $d.Range($bmStart, $bmStart + 9)
    $essentialRange.Font.Color = $RED

    $newBmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $newBmRange)
}
